$wb = $excel.ActiveWorkbook

# --- Metrics sheet: update the raw input values in column B ---
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value = 120285.72
$metrics.Range("B3").Value = 103388.09000000001
$metrics.Range("B4").Value = 36992.519999999997
$metrics.Range("B5").Value = 4899
$metrics.Range("B6").Value = 5322992.830000001
$metrics.Range("B7").Value = 4503741.0500000007
$metrics.Range("B8").Value = 1568949.4000000004
$metrics.Range("B9").Value = 207606
$metrics.Range("B10").Value = 33788373.819999993
$metrics.Range("B11").Value = 31779016.210000001
$metrics.Range("B12").Value = 11850671.439999996
$metrics.Range("B13").Value = 1305236

# Move the active selection on the Metrics sheet to D9
$metrics.Activate()
$metrics.Range("D9").Select()

# --- today sheet: move the active selection to E6 ---
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("E6").Select()

$wb.Save()
